$wb = $excel.ActiveWorkbook

# Sheet 1: "VENTAS POR GRUPO"
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("I4").Value = 129.36
$ws1.Range("I19").Value = "1 de 17"

# Sheet 2: "VENTA MENSUAL"
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 877.88
$ws2.Range("F19").Value = 22279.55

# Sheet 3: "CUMPLIMIENTO MENSUAL"
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D8").Value = 129.36
$ws3.Range("E8").Value = 495.64
$ws3.Range("F8").Value = 0.206976
$ws3.Range("D19").Value = 22279.55
$ws3.Range("E19").Value = 24939.75386304603
$ws3.Range("F19").Value = 0.4718313947325267
